# "Generate Report for Handback"
# The localization-status workbook gets its per-language sheets (zh-cn, de-de)
# updated with the handback results: the overall status text moves from
# "Ready for handoff" to "Handed back: in sync with en-US", the new
# "Latest Target File" / "Latest Handback File" columns (E/F) are populated
# with hyperlinked file names, and "Latest Handback DateTime" (G) is
# stamped with the real handback timestamp.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$langSheets = @(
    @{
        Name         = "zh-cn"
        HandbackTime = "2016-03-04 08:24:06"
        SourceUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/d8258c0f402b27f9a3e3adff1e69e6eb3e47ff5c/e2e/a.md"
        SourceUrlB   = "https://github.com/OpenLocalizationTest/oltest/blob/d8258c0f402b27f9a3e3adff1e69e6eb3e47ff5c/e2e/b.md"
        TargetDisp   = "a.md"
        HandoffUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/305becb69a9522d98652a60b98d10b3d515fc583/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
        HandoffDisp  = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
        ConfigUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/d8258c0f402b27f9a3e3adff1e69e6eb3e47ff5c/.localization-config"
        ConfigDisp   = ".localization-config"
    },
    @{
        Name         = "de-de"
        HandbackTime = "2016-03-04 08:24:24"
        SourceUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/d8258c0f402b27f9a3e3adff1e69e6eb3e47ff5c/e2e/a.md"
        SourceUrlB   = "https://github.com/OpenLocalizationTest/oltest/blob/d8258c0f402b27f9a3e3adff1e69e6eb3e47ff5c/e2e/b.md"
        TargetDisp   = "a.md"
        HandoffUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b57612b6101b810eee7a828a73f6844006dada43/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
        HandoffDisp  = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
        ConfigUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/d8258c0f402b27f9a3e3adff1e69e6eb3e47ff5c/.localization-config"
        ConfigDisp   = ".localization-config"
    }
)

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)

    # Populate the new "Latest Target File" (E) / "Latest Handback File" (F)
    # columns and stamp "Latest Handback DateTime" (G) now that handback
    # happened.
    $ws.Range("G2").Value = $lang.HandbackTime
    $ws.Range("G3").Value = $lang.HandbackTime

    # Status column (B) reflects the handback instead of "ready for handoff".
    $ws.Range("B2:B3").Value = $newStatus

    # Rebuild the hyperlinks collection so the new cells slot in next to
    # their row (A2,C2,E2,F2,A3,C3,E3,F3,A4) instead of being tacked on at
    # the end.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $lang.SourceUrl, "", "", "a.md")
    $ws.Hyperlinks.Add($ws.Range("C2"), $lang.HandoffUrl, "", "", $lang.HandoffDisp)
    $ws.Hyperlinks.Add($ws.Range("E2"), $lang.SourceUrl, "", "", $lang.TargetDisp)
    $ws.Hyperlinks.Add($ws.Range("F2"), $lang.HandoffUrl, "", "", $lang.HandoffDisp)
    $ws.Hyperlinks.Add($ws.Range("A3"), $lang.SourceUrlB, "", "", "b.md")
    $ws.Hyperlinks.Add($ws.Range("C3"), $lang.HandoffUrl, "", "", $lang.HandoffDisp)
    $ws.Hyperlinks.Add($ws.Range("E3"), $lang.SourceUrl, "", "", $lang.TargetDisp)
    $ws.Hyperlinks.Add($ws.Range("F3"), $lang.HandoffUrl, "", "", $lang.HandoffDisp)
    $ws.Hyperlinks.Add($ws.Range("A4"), $lang.ConfigUrl, "", "", $lang.ConfigDisp)
}
